$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting -------------------------------------------------------
# Rows 6-7 reuse the existing "date, centered" style that's already used
# by rows 2-5 (column E already carries it) -- copy format only so the
# existing style index is reused instead of a new one being created.
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E6:G7").PasteSpecial(-4122) | Out-Null

# Rows 8-12 get a new style: date format, no special alignment. Apply the
# number format to a single cell first (this is what creates the single
# new style record), then clone that exact formatting onto the remaining
# cells via copy/paste-special so no further styles get created.
$ws.Range("E8").NumberFormat = "mm-dd-yy"
$ws.Range("E8").Copy() | Out-Null
$ws.Range("E8:G11").PasteSpecial(-4122) | Out-Null
$ws.Range("E12:F12").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Values -------------------------------------------------------------
$ws.Range("E6").Value = 44496
$ws.Range("F6").Value = 44502
$ws.Range("G6").Value = 44502

$ws.Range("E7").Value = 44502
$ws.Range("F7").Value = 44507
$ws.Range("G7").Value = 44507

$ws.Range("E8").Value = 44508
$ws.Range("F8").Value = 44512
$ws.Range("G8").Value = 44512

$ws.Range("E9").Value = 44513
$ws.Range("F9").Value = 44517
$ws.Range("G9").Value = 44517

$ws.Range("E10").Value = 44518
$ws.Range("F10").Value = 44523
$ws.Range("G10").Value = 44523

$ws.Range("E11").Value = 44524
$ws.Range("F11").Value = 44529
$ws.Range("G11").Value = 44529

$ws.Range("E12").Value = 44532
$ws.Range("F12").Value = 44542

# --- Selection ------------------------------------------------------------
$ws.Range("G12").Select() | Out-Null
